$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("RegisterUserTest")

# Capture the row-3 values ("EditUser" PUT request) before it is removed from
# RegisterUserTest, so they can be re-homed on the new sheet.
$a3 = $ws2.Range("A3").Value2
$b3 = $ws2.Range("B3").Value2
$c3 = $ws2.Range("C3").Value2
$d3 = $ws2.Range("D3").Value2

# Also grab the shared header row text so the new sheet mirrors it exactly.
$h1 = $ws2.Range("A1").Value2
$h2 = $ws2.Range("B1").Value2
$h3 = $ws2.Range("C1").Value2
$h4 = $ws2.Range("D1").Value2

# Capture page setup (margins/header/footer) so the new sheet matches the
# look of its siblings instead of falling back to engine defaults.
$pmLeft = $ws2.PageSetup.LeftMargin
$pmRight = $ws2.PageSetup.RightMargin
$pmTop = $ws2.PageSetup.TopMargin
$pmBottom = $ws2.PageSetup.BottomMargin
$pmHeader = $ws2.PageSetup.HeaderMargin
$pmFooter = $ws2.PageSetup.FooterMargin
$cHeader = $ws2.PageSetup.CenterHeader
$cFooter = $ws2.PageSetup.CenterFooter

# Add the new worksheet as the last tab.
$newWs = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$newWs.Name = "EditUserTest"

$newWs.PageSetup.LeftMargin = $pmLeft
$newWs.PageSetup.RightMargin = $pmRight
$newWs.PageSetup.TopMargin = $pmTop
$newWs.PageSetup.BottomMargin = $pmBottom
$newWs.PageSetup.HeaderMargin = $pmHeader
$newWs.PageSetup.FooterMargin = $pmFooter
$newWs.PageSetup.CenterHeader = $cHeader
$newWs.PageSetup.CenterFooter = $cFooter

$newWs.Range("A1").Value2 = $h1
$newWs.Range("B1").Value2 = $h2
$newWs.Range("C1").Value2 = $h3
$newWs.Range("D1").Value2 = $h4

$newWs.Range("A2").Value2 = $a3
$newWs.Range("B2").Value2 = $b3
$newWs.Range("C2").Value2 = $c3
$newWs.Range("D2").Value2 = $d3

$newWs.Range("B2").Select() | Out-Null

# Remove the row from RegisterUserTest now that it lives on its own sheet,
# and correct the response code for the remaining row.
$ws2.Rows("3:3").Delete()
$ws2.Range("C2").Value2 = 200

# Keep RegisterUserTest the active/selected tab, matching the original
# workbook state, with the cursor parked on the edited cell.
$ws2.Activate()
$ws2.Range("C2").Select() | Out-Null
